# Clear the "*" placeholder values out of the grip_open / grip_closed rows
# (C11:G11 and C12:G12). These cells used to hold the literal string "*" as
# a "don't care" sentinel for the s0-s4 servo columns; the consuming code
# now does int conversion + NaN checking instead, so the cells should be
# truly empty rather than containing a stray string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("C11:G12")
$range.ClearContents()
$range.NumberFormat = "General"
